$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data (e.g. "65.901.61"),
# including values that look like ordinary decimal numbers (e.g. "581.39"). Excel would
# normally auto-convert a plain numeric-looking string typed into a General-formatted cell
# into a real number, which would lose the original text representation (e.g. drop trailing
# zeros like "159.00" -> 159). Mark those specific cells as Text first so the new price
# strings are stored verbatim, matching the existing inline-string cells.
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '65.901.61'
$ws.Cells.Item(2, 5).Value = '  -2.56%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.463.38'

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '581.39'
$ws.Cells.Item(5, 5).Value = '  -1.58%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '172.76'
$ws.Cells.Item(6, 5).Value = '  -3.36%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -2.26%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '3.464.36'
$ws.Cells.Item(9, 5).Value = '  +0.67%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -5.93%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '6.85'
$ws.Cells.Item(11, 5).Value = '  -1.73%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -4.10%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '4.060.09'
$ws.Cells.Item(13, 5).Value = '  +0.52%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.84%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '29.75'
$ws.Cells.Item(15, 5).Value = '  -7.10%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '65.980.14'
$ws.Cells.Item(16, 5).Value = '  -2.40%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -3.26%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.461.22'
$ws.Cells.Item(18, 5).Value = '  +0.56%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -3.58%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.17%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '365.53'
$ws.Cells.Item(21, 5).Value = '  -6.07%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '7.69'
$ws.Cells.Item(22, 5).Value = '  -1.85%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.08%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '72.04'
$ws.Cells.Item(24, 5).Value = '  +0.98%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '0.532'
$ws.Cells.Item(25, 5).Value = '  -0.02%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +3.73%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '9.52'
$ws.Cells.Item(27, 5).Value = '  -7.02%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.19%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.00%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '23.79'
$ws.Cells.Item(30, 5).Value = '  +2.48%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '5.74'
$ws.Cells.Item(31, 5).Value = '  -4.93%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -3.22%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.00%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Fetch.AI'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(34, 4).Value = '1.29'
$ws.Cells.Item(34, 5).Value = '  -6.87%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Aptos'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(35, 4).Value = '7.08'
$ws.Cells.Item(35, 5).Value = '  -1.81%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.91%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '159.00'
$ws.Cells.Item(37, 5).Value = '  -1.28%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'EnergySwap'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(38, 4).Value = '29.03'
$ws.Cells.Item(38, 5).Value = '  +12.68%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Mantle'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(39, 4).Value = '0.886'
$ws.Cells.Item(39, 5).Value = '  +0.40%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '2.795.59'
$ws.Cells.Item(40, 5).Value = '  +3.68%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '1.76'
$ws.Cells.Item(41, 5).Value = '  -5.30%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(42, 4).Value = '2.55'
$ws.Cells.Item(42, 5).Value = '  -6.72%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).Value = '6.43'
$ws.Cells.Item(43, 5).Value = '  -2.49%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -3.62%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.0679'
$ws.Cells.Item(45, 5).Value = '  -4.70%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -2.80%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '24.06'
$ws.Cells.Item(47, 5).Value = '  -7.34%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -2.85%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '311.12'
$ws.Cells.Item(49, 5).Value = '  -3.96%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -2.57%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -3.19%  '
